$d = $word.ActiveDocument

# The whole SmartRewards FAQ block (everything between the "SMARTREWARDS
# FAQ'S" heading and the trailing blank paragraph) is being removed. Locate
# it by its first and last paragraph's text so the script does not depend on
# hard-coded paragraph indices.
$startMarker = "How much SmartRewards can I get each month?"
$endMarker = "SmartRewards will then be distributed after every 47500 Blocks to all eligible addresses. The payouts will occur 200 Blocks after the cycle ends and every second block 1000 addresses will get paid."

$startIdx = -1
$endIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($startIdx -eq -1 -and $t -like "*$startMarker*") {
        $startIdx = $i
    }
    if ($t -like "*$endMarker*") {
        $endIdx = $i
    }
}

if ($startIdx -ne -1 -and $endIdx -ne -1) {
    $deleteStart = $d.Paragraphs.Item($startIdx).Range.Start
    $deleteEnd = $d.Paragraphs.Item($endIdx).Range.End
    $deleteRange = $d.Range($deleteStart, $deleteEnd)
    $deleteRange.Delete()
}
